# fix bug siswa dan rombel
#
# The "telp siswa" column (AA) was a mistake: the real data that gets
# collected for a student is telp ayah / telp ibu / telp wali, not a
# separate "telp siswa" entry, so every header from AA onwards shifts
# left by one and a new "telp wali" column is appended at the end (AC).
#
#   AA: telp siswa -> telp ayah
#   AB: telp ayah   -> telp ibu
#   AC: telp ibu    -> telp wali

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row relabeling -------------------------------------------------
$ws.Range("AA1").Value = "telp ayah"
$ws.Range("AB1").Value = "telp ibu"
$ws.Range("AC1").Value = "telp wali"

# --- keep the column input-message hints (data validation prompts) in sync -
$ws.Range("AA1:AA1048576").Validation.InputMessage = "telp ayah"
$ws.Range("AB1:AB1048576").Validation.InputMessage = "telp ibu"
$ws.Range("AC1:AC1048576").Validation.InputMessage = "telp wali"

# --- restore the view/selection the author left the sheet in ---------------
# (scrolled right so column M is at the left edge, with AA2 selected)
$excel.ActiveWindow.ScrollColumn = 13
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AA2").Select() | Out-Null
